# Production-Follow-up workbook update:
#  - Insert a new title row at the top of the main table and label it
#    "Plan Vs Achievement"
#  - Everything that used to live in rows 1-21 shifts down to rows 2-22
#  - Restore (best-effort) the view state (scroll position / selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Production-Follow-up")
$ws.Activate()

# Insert a new blank row above the current row 1; this pushes every
# existing row (data table, forecast table, etc.) down by one, which is
# exactly what the diff shows (row N -> row N+1 everywhere).
$ws.Rows.Item(1).Insert()

# Give the newly-inserted row 1 its label.
$ws.Range("A1").Value = "Plan Vs Achievement"

# Best-effort restore of the saved view state (scrolled so column E is at
# the left edge, with S4 as the active cell of the remembered selection).
$ws.Range("S6:S12").Select()
$r1 = $ws.Range("S6:S12")
$r2 = $ws.Range("S4")
$u = $excel.Union($r1, $r2)
$u.Select()
